# "Problems with big memory"
#
# The Configuration sheet's "Categorized Columns" rule table used
# auto-generated, meaningless section headers ("Rules - <row>") and
# encoded each rule's output as an opaque row-derived number. This
# replaces the headers with the real source-column name and replaces
# the numeric outputs with descriptive category labels, and refreshes
# the Walk_times / Walk_mins rule blocks (moved further down the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# --- Rename generic "Rules - N" section headers to the real column name,
#     and replace the numeric rule-output codes with descriptive labels ---
$ws.Range("B10").Value = "Rules - Age"
$ws.Range("D11").Value = "Teenage"
$ws.Range("D12").Value = "20s"
$ws.Range("D13").Value = "30s"
$ws.Range("D14").Value = "old"

$ws.Range("B15").Value = "Rules - Sex"
$ws.Range("D16").Value = "Male"
$ws.Range("D17").Value = "Female"

$ws.Range("B18").Value = "Rules - BMI"
$ws.Range("D19").Value = "Underweight"
$ws.Range("D20").Value = "Healthy"
$ws.Range("D21").Value = "Overweight"
$ws.Range("D22").Value = "Obese"

$ws.Range("B23").Value = "Rules - Marital"
$ws.Range("D24").Value = "Never"
$ws.Range("D25").Value = "Married"
$ws.Range("D26").Value = "Once"

$ws.Range("B27").Value = "Rules - Children"
$ws.Range("D28").Value = "None"
$ws.Range("D29").Value = "Yes"

$ws.Range("B30").Value = "Rules - Sleep"
$ws.Range("D31").Value = "Poor"
$ws.Range("D32").Value = "Enough"

$ws.Range("B33").Value = "Rules - Pt_ft"
$ws.Range("D34").Value = "Part"
$ws.Range("D35").Value = "Full"

$ws.Range("B36").Value = "Rules - Cigs"
$ws.Range("D37").Value = "No"
$ws.Range("D38").Value = "Yes"

$ws.Range("B39").Value = "Rules - Chol"
$ws.Range("D40").Value = "Yes"
$ws.Range("D41").Value = "No"

$ws.Range("B42").Value = "Rules - Fruit"
$ws.Range("D43").Value = "No"
$ws.Range("D44").Value = "One"
$ws.Range("D45").Value = "Many"

$ws.Range("B46").Value = "Rules - Veg"
$ws.Range("D47").Value = "No"
$ws.Range("D48").Value = "One"
$ws.Range("D49").Value = "Many"

# --- The old Walk_times (B55:D59) / Walk_mins (B50:D53) rule blocks are
#     cleared out of their old spot ... ---
$ws.Range("B50").ClearContents()
$ws.Range("C51").ClearContents()
$ws.Range("D51").ClearContents()
$ws.Range("C52").ClearContents()
$ws.Range("D52").ClearContents()
$ws.Range("C53").ClearContents()
$ws.Range("D53").ClearContents()

$ws.Range("B55").ClearContents()
$ws.Range("C56").ClearContents()
$ws.Range("D56").ClearContents()
$ws.Range("C57").ClearContents()
$ws.Range("D57").ClearContents()
$ws.Range("C58").ClearContents()
$ws.Range("D58").ClearContents()
$ws.Range("C59").ClearContents()
$ws.Range("D59").ClearContents()

# ... and re-entered further down (rows 60-68), Walk_times first, then
# Walk_mins, each now carrying descriptive category labels.
$ws.Range("B60").Value = "Rules - Walk_times"
$ws.Range("C61").Value = "value == 0"
$ws.Range("D61").Value = "Zero"
$ws.Range("C62").Value = "value <= 7"
$ws.Range("D62").Value = "Commuter"
$ws.Range("C63").Value = "value > 7 && value <= 14"
$ws.Range("D63").Value = "Walk"
$ws.Range("C64").Value = "value > 14"
$ws.Range("D64").Value = "Walker"

$ws.Range("B65").Value = "Rules - Walk_mins"
$ws.Range("C66").Value = "value <= 70"
$ws.Range("D66").Value = "CouchPotato"
$ws.Range("C67").Value = "value > 70 && value <= 140"
$ws.Range("D67").Value = "Regular"
$ws.Range("C68").Value = "value > 140"
$ws.Range("D68").Value = "Enough"

# A couple of section headers got a touch more breathing room.
$ws.Rows.Item(36).RowHeight = 10.85
$ws.Rows.Item(39).RowHeight = 10.85
$ws.Rows.Item(65).RowHeight = 10.85
$ws.Rows.Item(66).RowHeight = 10.85

# Leave the sheet scrolled to the top with the new block selected.
$ws.Activate()
$ws.Range("A49").Select()
